# Applies the row additions / cell removals described in the commit:
#   - row 152: AA152 / AD152 (empty placeholder cells) are dropped
#   - rows 153-155 are appended with sensor-log data (all text cells,
#     matching the workbook's existing "everything is text" convention)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 152 previously had empty placeholder cells in AA/AD; overwriting drops them.
$ws.Range("AA152").ClearContents()
$ws.Range("AD152").ClearContents()

# Row 153
$row153 = [ordered]@{
    'A' = '2022-06-16 08:54:37'
    'B' = '21.1'
    'C' = '25.8'
    'D' = '20.0'
    'E' = '.'
    'F' = '24.5'
    'G' = '.'
    'H' = '46'
    'I' = '.'
    'J' = '50'
    'K' = '50'
    'L' = '0'
    'M' = '0'
    'N' = '50'
    'O' = '50'
    'P' = '.'
    'Q' = '.'
    'R' = '0'
    'S' = '0'
    'T' = '0'
    'U' = '0'
    'V' = '100'
    'W' = '51'
    'X' = '.'
    'Y' = '.'
    'Z' = '10.28'
    'AB' = '7002200'
    'AC' = '2200'
    'AE' = '21.2'
    'AF' = '25.8'
    'AG' = '20.0'
    'AH' = '50'
    'AI' = '50'
    'AJ' = '50'
    'AK' = '50'
    'AL' = '50'
    'AM' = '51'
    'AN' = '.'
    'AO' = '0'
    'AP' = '26'
    'AQ' = '0'
    'AR' = '0.16'
    'AS' = '0.15'
    'AT' = '0.12'
    'AU' = '33.25'
    'AV' = '3568.69'
    'AW' = '0.00'
    'AX' = '17.78'
    'AY' = '2913.03'
    'AZ' = '0.00'
    'BA' = '87.68'
    'BB' = '10317.52'
    'BC' = '19.0'
    'BD' = '.'
    'BE' = '.'
    'BF' = '46'
    'BG' = '.'
    'BH' = '46'
    'BI' = '146030596'
}
foreach ($col in $row153.Keys) {
    $cell = $ws.Range($col + "153")
    $cell.NumberFormat = "@"
    $cell.Value = $row153[$col]
}

# Row 154
$row154 = [ordered]@{
    'A' = '2022-06-16 08:56:16'
    'B' = '21.3'
    'C' = '25.9'
    'D' = '20.2'
    'E' = '.'
    'F' = '24.5'
    'G' = '.'
    'H' = '46'
    'I' = '.'
    'J' = '50'
    'K' = '50'
    'L' = '0'
    'M' = '0'
    'N' = '50'
    'O' = '50'
    'P' = '.'
    'Q' = '.'
    'R' = '0'
    'S' = '0'
    'T' = '0'
    'U' = '0'
    'V' = '100'
    'W' = '51'
    'X' = '.'
    'Y' = '.'
    'Z' = '10.28'
    'AB' = '7002200'
    'AC' = '2200'
    'AE' = '21.4'
    'AF' = '25.9'
    'AG' = '20.2'
    'AH' = '50'
    'AI' = '50'
    'AJ' = '50'
    'AK' = '50'
    'AL' = '50'
    'AM' = '51'
    'AN' = '.'
    'AO' = '0'
    'AP' = '26'
    'AQ' = '0'
    'AR' = '0.16'
    'AS' = '0.15'
    'AT' = '0.12'
    'AU' = '33.25'
    'AV' = '3568.69'
    'AW' = '0.00'
    'AX' = '17.78'
    'AY' = '2913.03'
    'AZ' = '0.00'
    'BA' = '87.68'
    'BB' = '10317.52'
    'BC' = '19.0'
    'BD' = '.'
    'BE' = '.'
    'BF' = '46'
    'BG' = '.'
    'BH' = '46'
    'BI' = '146030596'
}
foreach ($col in $row154.Keys) {
    $cell = $ws.Range($col + "154")
    $cell.NumberFormat = "@"
    $cell.Value = $row154[$col]
}

# Row 155
$row155 = [ordered]@{
    'A' = '2022-06-16 18:19:22'
    'B' = '23.4'
    'C' = '26.4'
    'D' = '22.1'
    'E' = '.'
    'F' = '24.9'
    'G' = '.'
    'H' = '43'
    'I' = '.'
    'J' = '50'
    'K' = '50'
    'L' = '0'
    'M' = '0'
    'N' = '50'
    'O' = '50'
    'P' = '.'
    'Q' = '.'
    'R' = '0'
    'S' = '0'
    'T' = '0'
    'U' = '0'
    'V' = '100'
    'W' = '51'
    'X' = '.'
    'Y' = '.'
    'Z' = '9.84'
    'AB' = '7002200'
    'AC' = '2200'
    'AE' = '23.3'
    'AF' = '26.4'
    'AG' = '22.2'
    'AH' = '50'
    'AI' = '50'
    'AJ' = '50'
    'AK' = '50'
    'AL' = '50'
    'AM' = '51'
    'AN' = '.'
    'AO' = '0'
    'AP' = '26'
    'AQ' = '0'
    'AR' = '0.16'
    'AS' = '0.16'
    'AT' = '0.36'
    'AU' = '33.50'
    'AV' = '3568.94'
    'AW' = '0.00'
    'AX' = '17.78'
    'AY' = '2913.03'
    'AZ' = '0.01'
    'BA' = '87.69'
    'BB' = '10317.53'
    'BC' = '19.0'
    'BD' = '.'
    'BE' = '.'
    'BF' = '43'
    'BG' = '.'
    'BH' = '43'
    'BI' = '146030596'
}
foreach ($col in $row155.Keys) {
    $cell = $ws.Range($col + "155")
    $cell.NumberFormat = "@"
    $cell.Value = $row155[$col]
}

